# Team-Meeting-11.docx regeneration:
#  1. Remove the three leading "Home" / "<- Back to Home" / "Download Word
#     Document" hyperlink paragraphs that used to sit at the very top of
#     the document, right before the "Team Meeting Agenda ..." Heading1
#     paragraph.
#  2. Regenerate table formatting: every w:tbl in the body should carry a
#     preferred width of 100% (OOXML <w:tblW w:type="pct" w:w="5000"/>)
#     instead of the "auto / 0" width that was there before.

$d = $word.ActiveDocument

# --- 1. Drop the hyperlink paragraphs that precede the first heading ---
$headingRange = $d.Content
$headingRange.Find.Execute("Team Meeting Agenda - Eleventh Team Meeting") | Out-Null

if ($headingRange.Start -gt 0) {
    $d.Range(0, $headingRange.Start).Delete()
}

# --- 2. Force every table to a 100% preferred width ---
# wdPreferredWidthPercent = 2; COM PreferredWidth is expressed in
# twentieths of a percent, so 250 -> OOXML w:w="5000" (100%).
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $tbl = $d.Tables($i)
    $tbl.PreferredWidthType = 2
    $tbl.PreferredWidth = 250
}

Write-Output ("Paragraphs now: " + $d.Paragraphs.Count)
Write-Output ("Tables updated: " + $d.Tables.Count)
